$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("F2").Value = 25.56000000000056
$ws.Range("H2").Value = 0.01898410260938044
$ws.Range("I2").Value = 0.01898410260938044
$ws.Range("L2").Value = 8.953648041754393
$ws.Range("M2").Value = "[1.1814557756461035, 16.725840307862683]"
$ws.Range("N2").Value = 0.02491860708663762
$ws.Range("O2").Value = 0.02491860708663762
$ws.Range("P2").Value = -1.773631888513156
$ws.Range("Q2").Value = "[-2.918316298830158, -0.6289474781961548]"
$ws.Range("R2").Value = 0.003146685633408808
$ws.Range("S2").Value = 0.003146685633408808
$ws.Range("T2").Value = 12.93432694348022
$ws.Range("U2").Value = "[8.50525445095578, 17.36339943600465]"
$ws.Range("V2").Value = [double]"4.674902871482089e-07"
$ws.Range("W2").Value = [double]"4.674902871482089e-07"
$ws.Range("X2").Value = 7.215135135135291
$ws.Range("Y2").Value = 2.558558558558614
$ws.Range("Z2").Value = 11.87171171171197
$ws.Range("F3").Value = 25.56000000000056
$ws.Range("H3").Value = 0.3788083684975772
$ws.Range("I3").Value = 0.3788083684975772
$ws.Range("L3").Value = 4.603460538935749
$ws.Range("M3").Value = "[-4.0938923272132985, 13.300813405084797]"
$ws.Range("N3").Value = 0.2920858009142597
$ws.Range("O3").Value = 0.2920858009142597
$ws.Range("P3").Value = -2.062947728483388
$ws.Range("Q3").Value = "[-5.1196324725167, 0.9937370155499243]"
$ws.Range("R3").Value = 0.1808209296953545
$ws.Range("S3").Value = 0.1808209296953545
$ws.Range("T3").Value = 11.08903307907023
$ws.Range("U3").Value = "[6.416397800184022, 15.761668357956445]"
$ws.Range("V3").Value = [double]"1.907523778976739e-05"
$ws.Range("W3").Value = [double]"1.907523778976739e-05"
$ws.Range("X3").Value = 8.392072072072256
$ws.Range("Y3").Value = -4.042522522522608
$ws.Range("Z3").Value = 20.82666666666712
$ws.Range("F4").Value = 25.56000000000056
$ws.Range("H4").Value = 0.4509490293006349
$ws.Range("I4").Value = 0.4509490293006349
$ws.Range("L4").Value = 3.312331369592845
$ws.Range("M4").Value = "[-3.1819920051896258, 9.806654744375315]"
$ws.Range("N4").Value = 0.30978700324876
$ws.Range("O4").Value = 0.30978700324876
$ws.Range("P4").Value = -2.842842601446619
$ws.Range("Q4").Value = "[-5.981290517645432, 0.2956053147521933]"
$ws.Range("R4").Value = 0.07473565414383554
$ws.Range("S4").Value = 0.07473565414383554
$ws.Range("T4").Value = 10.45272191066313
$ws.Range("U4").Value = "[6.867311229718471, 14.038132591607788]"
$ws.Range("V4").Value = [double]"4.837552800562861e-07"
$ws.Range("W4").Value = [double]"4.837552800562861e-07"
$ws.Range("X4").Value = 11.56468468468493
$ws.Range("Y4").Value = -1.202522522522552
$ws.Range("Z4").Value = 24.33189189189242
$ws.Range("F5").Value = 25.56000000000056
$ws.Range("H5").Value = 0.1040681842906682
$ws.Range("I5").Value = 0.1040681842906682
$ws.Range("L5").Value = 6.329898694252444
$ws.Range("M5").Value = "[-1.3424619929676052, 14.002259381472493]"
$ws.Range("N5").Value = 0.1035273658873042
$ws.Range("O5").Value = 0.1035273658873042
$ws.Range("P5").Value = -2.528368862348542
$ws.Range("Q5").Value = "[-5.377500938577123, 0.3207632138800385]"
$ws.Range("R5").Value = 0.08061708650077115
$ws.Range("S5").Value = 0.08061708650077115
$ws.Range("T5").Value = 10.86734419699205
$ws.Range("U5").Value = "[6.837092904602738, 14.897595489381363]"
$ws.Range("V5").Value = [double]"2.164761783607361e-06"
$ws.Range("W5").Value = [double]"2.164761783607361e-06"
$ws.Range("X5").Value = 10.28540540540563
$ws.Range("Y5").Value = -1.304864864864891
$ws.Range("Z5").Value = 21.87567567567615
$ws.Range("H6").Value = 0.03320269698579303
$ws.Range("I6").Value = 0.03320269698579303
$ws.Range("L6").Value = 8.191444379212301
$ws.Range("M6").Value = "[-0.13949730943580718, 16.52238606786041]"
$ws.Range("N6").Value = 0.05379458074699062
$ws.Range("O6").Value = 0.05379458074699062
$ws.Range("P6").Value = -2.465474114528927
$ws.Range("Q6").Value = "[-4.276842851733853, -0.6541053773240018]"
$ws.Range("R6").Value = 0.008745071570594343
$ws.Range("S6").Value = 0.008745071570594343
$ws.Range("T6").Value = 12.18865981368737
$ws.Range("U6").Value = "[7.8832926343404806, 16.494026993034254]"
$ws.Range("V6").Value = [double]"8.628974510482834e-07"
$ws.Range("W6").Value = [double]"8.628974510482834e-07"
$ws.Range("X6").Value = 9.021101101101165
$ws.Range("Y6").Value = 2.393353353353374
$ws.Range("Z6").Value = 15.64884884884896
$ws.Range("H7").Value = 0.2906456528647204
$ws.Range("I7").Value = 0.2906456528647204
$ws.Range("L7").Value = 4.72178380957028
$ws.Range("M7").Value = "[-2.7270461741641423, 12.170613793304703]"
$ws.Range("N7").Value = 0.2082461008887173
$ws.Range("O7").Value = 0.2082461008887173
$ws.Range("P7").Value = -3.069263693597235
$ws.Range("Q7").Value = "[-6.188843185450164, 0.05031579825569299]"
$ws.Range("R7").Value = 0.05365082899869034
$ws.Range("S7").Value = 0.05365082899869034
$ws.Range("T7").Value = 10.11293356169957
$ws.Range("U7").Value = "[5.814290085741356, 14.411577037657786]"
$ws.Range("V7").Value = [double]"2.186551073957155e-05"
$ws.Range("W7").Value = [double]"2.186551073957155e-05"
$ws.Range("X7").Value = 11.23035035035043
$ws.Range("Y7").Value = -0.184104104104108
$ws.Range("Z7").Value = 22.64480480480496
$ws.Range("H8").Value = 0.06207026729265031
$ws.Range("I8").Value = 0.06207026729265031
$ws.Range("L8").Value = 7.8307154748271
$ws.Range("M8").Value = "[-0.9674517829569567, 16.628882732611157]"
$ws.Range("N8").Value = 0.07975562540270054
$ws.Range("O8").Value = 0.07975562540270054
$ws.Range("P8").Value = -2.176158274558696
$ws.Range("Q8").Value = "[-4.415211296937008, 0.06289474781961601]"
$ws.Range("R8").Value = 0.05650466943310617
$ws.Range("S8").Value = 0.05650466943310617
$ws.Range("T8").Value = 11.27025319171193
$ws.Range("U8").Value = "[6.733392893584361, 15.807113489839498]"
$ws.Range("V8").Value = [double]"9.097349810360811e-06"
$ws.Range("W8").Value = [double]"9.097349810360811e-06"
$ws.Range("X8").Value = 7.962502502502558
$ws.Range("Y8").Value = -0.2301301301301333
$ws.Range("Z8").Value = 16.15513513513525
$ws.Range("H9").Value = 0.2973574546816372
$ws.Range("I9").Value = 0.2973574546816372
$ws.Range("L9").Value = 4.457914812195375
$ws.Range("M9").Value = "[-2.7134458715515954, 11.629275495942345]"
$ws.Range("N9").Value = 0.2170304856743204
$ws.Range("O9").Value = 0.2170304856743204
$ws.Range("P9").Value = -1.761052938949233
$ws.Range("Q9").Value = "[-4.698237662125276, 1.1761317842268095]"
$ws.Range("R9").Value = 0.2335102063558581
$ws.Range("S9").Value = 0.2335102063558581
$ws.Range("T9").Value = 11.75257518760518
$ws.Range("U9").Value = "[7.8069592984286205, 15.698191076781747]"
$ws.Range("V9").Value = [double]"3.129946348767731e-07"
$ws.Range("W9").Value = [double]"3.129946348767731e-07"
$ws.Range("X9").Value = 6.443643643643686
$ws.Range("Y9").Value = -4.303433433433463
$ws.Range("Z9").Value = 17.19072072072083
$ws.Range("H10").Value = 0.08222511396056253
$ws.Range("I10").Value = 0.08222511396056253
$ws.Range("L10").Value = 7.62723278651807
$ws.Range("M10").Value = "[-1.0879717287343178, 16.34243730177046]"
$ws.Range("N10").Value = 0.08474679779399419
$ws.Range("O10").Value = 0.08474679779399419
$ws.Range("P10").Value = -1.270473905956233
$ws.Range("Q10").Value = "[-3.00636894577762, 0.4654211338651546]"
$ws.Range("R10").Value = 0.1474191531253826
$ws.Range("S10").Value = 0.1474191531253826
$ws.Range("T10").Value = 12.73685706324309
$ws.Range("U10").Value = "[8.079941435547504, 17.393772690938675]"
$ws.Range("V10").Value = [double]"1.663891038283793e-06"
$ws.Range("W10").Value = [double]"1.663891038283793e-06"
$ws.Range("X10").Value = 4.648628628628661
$ws.Range("Y10").Value = -1.702962962962974
$ws.Range("Z10").Value = 11.0002202202203
